# Add a new column E ("Variável alternativa") to sheet1, with a cyclic
# A/B/C/D value for every data row (rows 2-57), mirroring column D's header
# style for the E1 header cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell E1, styled like the other header cells (A1:D1)
$ws.Range("E1").Value = "Variável alternativa"
$ws.Range("E1").Font.Bold = $true
$ws.Range("E1").HorizontalAlignment = -4108  # xlCenter

# Cyclic A/B/C/D values for data rows 2 through 57
$labels = @("A", "B", "C", "D")
for ($row = 2; $row -le 57; $row++) {
    $label = $labels[($row - 2) % 4]
    $ws.Cells.Item($row, 5).Value = $label
}
